# Apply the commit's changes:
#  - Sheet1!E2 ("count" for "loginlogout" row) changes from 2 to 1
#    (kept as quoted/text "1", matching its existing quote-prefixed style)
#  - The active tab moves from "DATA" (sheet index 1) to "Sheet1" (index 0)
#  - The DATA sheet's selection changes from C12 to A2:F2 (active cell A2)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("DATA")

# Sheet1: E2 "2" -> "1" (leading apostrophe keeps it text-typed, same as before)
$ws1.Range("E2").Formula = "'1"

# DATA sheet: update its selection while it is still the active sheet so the
# selection sticks to that sheet's view state.
$ws2.Activate()
[void]$ws2.Range("A2:F2").Select()

# Finally make Sheet1 the active tab (tabSelected on Sheet1, removed from DATA).
$ws1.Activate()
